$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh aggregate stats after closing trade #29
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.25   # Current Capital
$summary.Range("B4").Value = -0.76     # Total P&L $
$summary.Range("B5").Value = -0.52     # Total P&L %
$summary.Range("B6").Value = 29        # Total Trades
$summary.Range("B8").Value = 17        # Losing Trades
$summary.Range("B9").Value = 34.48     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.25      # Capital
$status.Range("D4").Value = 29         # Trades
$status.Range("E4").Value = -0.76      # P&L $
$status.Range("F4").Value = -0.75      # P&L %
$status.Range("G4").Value = 34.48      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#29) to both "All Trades" and
# "MarketMaking" sheets (row 30).
# ---------------------------------------------------------------------------
$newTradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $newTradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(30, 1).Value = 29

    # Force the date-looking string to stay a plain text value instead of
    # being auto-converted into a date serial number.
    $ws.Cells.Item(30, 2).NumberFormat = "@"
    $ws.Cells.Item(30, 2).Value = "2026-02-17"
    $ws.Cells.Item(30, 2).Style = "Normal"

    $ws.Cells.Item(30, 3).Value = "13:19:16"
    $ws.Cells.Item(30, 4).Value = "MarketMaking"
    $ws.Cells.Item(30, 5).Value = "UP"
    $ws.Cells.Item(30, 6).Value = 0.08
    $ws.Cells.Item(30, 7).Value = 0.060995
    $ws.Cells.Item(30, 8).Value = "CLOSED"
    $ws.Cells.Item(30, 9).Value = -23.7566
    $ws.Cells.Item(30, 10).Value = -0.02
    $ws.Cells.Item(30, 11).Value = 99.25
    $ws.Cells.Item(30, 12).Value = 0
    $ws.Cells.Item(30, 13).Value = 0
    $ws.Cells.Item(30, 14).Value = 0.6
    $ws.Cells.Item(30, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(30, 16).Value = "early_exit"
    $ws.Cells.Item(30, 17).Value = 0.13
}
